# Fill in the "S3" (column D) scores for the three roster sheets that were
# still missing them (Astronauta, Senador, Mago). The "Resultados" column
# (L) already holds a shared IFERROR/SUM/COUNT formula over B:K, so it
# recalculates automatically once column D has values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Astronauta")
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("D5").Value = 1
$ws.Range("D7").Value = 0
$ws.Range("D8").Value = 1
$ws.Range("D9").Value = 1
$ws.Range("D11").Value = 1
$ws.Range("D12").Value = 0.5
$ws.Range("D13").Value = 0.5
$ws.Range("D14").Value = 0.5
$ws.Range("D15").Value = 0
$ws.Range("D16").Value = 0
$ws.Range("D17").Value = 1
$ws.Range("D18").Value = 0.6
$ws.Range("D19").Value = 0.4
$ws.Range("D20").Value = 1
$ws.Range("D21").Value = 0
$ws.Range("D22").Value = 0.8
$ws.Range("D23").Value = 1
$ws.Range("D24").Value = 1
$ws.Range("D25").Value = 0.6
$ws.Range("D26").Value = 1
$ws.Range("D27").Value = 0.7
$ws.Range("D28").Value = 0.7
$ws.Range("D29").Value = 0.7

$ws = $wb.Worksheets.Item("Senador")
$ws.Range("D2").Value = 0.5
$ws.Range("D3").Value = 0.5
$ws.Range("D5").Value = 1
$ws.Range("D7").Value = 0.7
$ws.Range("D8").Value = 0.7
$ws.Range("D9").Value = 1
$ws.Range("D11").Value = 1
$ws.Range("D12").Value = 1
$ws.Range("D13").Value = 1
$ws.Range("D14").Value = 1
$ws.Range("D15").Value = 0.25
$ws.Range("D16").Value = 1
$ws.Range("D17").Value = 1
$ws.Range("D18").Value = 1
$ws.Range("D19").Value = 0.5
$ws.Range("D20").Value = 1
$ws.Range("D22").Value = 1
$ws.Range("D23").Value = 1
$ws.Range("D24").Value = 1
$ws.Range("D25").Value = 0.3
$ws.Range("D26").Value = 0.6
$ws.Range("D27").Value = 1
$ws.Range("D28").Value = 0.7
$ws.Range("D29").Value = 0.8
$ws.Range("D30").Select()

$ws = $wb.Worksheets.Item("Mago")
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 0.5
$ws.Range("D5").Value = 1
$ws.Range("D7").Value = 0.5
$ws.Range("D8").Value = 0.9
$ws.Range("D9").Value = 1
$ws.Range("D11").Value = 1
$ws.Range("D12").Value = 0.5
$ws.Range("D13").Value = 0.5
$ws.Range("D14").Value = 1
$ws.Range("D15").Value = 0.3
$ws.Range("D16").Value = 0.6
$ws.Range("D17").Value = 1
$ws.Range("D18").Value = 1
$ws.Range("D19").Value = 0.5
$ws.Range("D20").Value = 1
$ws.Range("D22").Value = 1
$ws.Range("D23").Value = 1
$ws.Range("D24").Value = 1
$ws.Range("D25").Value = 0.5
$ws.Range("D26").Value = 0.2
$ws.Range("D27").Value = 0.6
$ws.Range("D28").Value = 0.6
$ws.Range("D29").Value = 0.7

# Reviewed "Mago" zoomed in, then moved on leaving the other sheets' and
# the workbook's selection/active-tab state as seen in the saved file.
$ws.Activate()
$excel.ActiveWindow.Zoom = 130
$ws.Range("D30").Select()

$ws = $wb.Worksheets.Item("Ninja")
$ws.Range("G2").Select()

$ws = $wb.Worksheets.Item("Astronauta")
$ws.Range("D30").Select()
